$wb = $excel.ActiveWorkbook

# Remember the sheet that was originally active so we can restore the
# active-tab state once we're done touching the "Modify Transaction" sheets.
$originalActiveSheet = $wb.ActiveSheet.Name

# --- "Modify Transaction" sheet (sheet4) ---------------------------------
$ws4 = $wb.Worksheets.Item("Modify Transaction")

# Insert a new row 2 (pushes the old "submitmakerepayment"/"click" row down
# to row 3) and populate it with the new ReceiptNumber scenario fields.
$ws4.Range("A2").EntireRow.Insert()
$ws4.Range("A2").Value = "ReceiptNumber"
$ws4.Range("B2").Value = 1234

# Update the recorded selection to match the new active cell.
[void]$ws4.Range("B10").Select()

# --- "Modify Transaction1" sheet (sheet5) --------------------------------
$ws5 = $wb.Worksheets.Item("Modify Transaction1")

$ws5.Range("A2").EntireRow.Insert()
$ws5.Range("A2").Value = "ReceiptNumber"
$ws5.Range("B2").Value = 432

[void]$ws5.Range("B7").Select()

# Restore the originally active sheet/tab.
[void]$wb.Worksheets.Item($originalActiveSheet).Select()
